$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 (was row 21) ---
$ws.Cells.Item(19,1).Value = 112390125
$ws.Cells.Item(19,2).Value = 56430
$ws.Cells.Item(19,3).Value = 'Ovaliderad'
$ws.Cells.Item(19,4).Value = 'NT'
$ws.Cells.Item(19,5).Value = 100109
$ws.Cells.Item(19,6).Value = 'Tretåig hackspett'
$ws.Cells.Item(19,7).Value = 'Picoides tridactylus'
$ws.Cells.Item(19,8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(19,13).Value = 'äldre spår'
$ws.Cells.Item(19,16).Value = 'Kärmsjöbäckens naturreservat (Kärmsjöbäckens naturreservat), Ång'
$ws.Cells.Item(19,17).Value = 583127
$ws.Cells.Item(19,18).Value = 7086576
$ws.Cells.Item(19,19).Value = 25
$ws.Cells.Item(19,20).Value = 'Västernorrland'
$ws.Cells.Item(19,21).Value = 'Sollefteå'
$ws.Cells.Item(19,22).Value = 'Ångermanland'
$ws.Cells.Item(19,23).Value = 'Junsele'
$ws.Cells.Item(19,25).NumberFormat = "@"
$ws.Cells.Item(19,25).Value = '2023-09-29'
$ws.Cells.Item(19,26).NumberFormat = "@"
$ws.Cells.Item(19,26).Value = '10:57'
$ws.Cells.Item(19,27).NumberFormat = "@"
$ws.Cells.Item(19,27).Value = '2023-09-29'
$ws.Cells.Item(19,28).NumberFormat = "@"
$ws.Cells.Item(19,28).Value = '10:57'
$ws.Cells.Item(19,30).Value = $false
$ws.Cells.Item(19,31).Value = $false
$ws.Cells.Item(19,33).Value = $false
$ws.Cells.Item(19,49).Value = 'Helena Thau'
$ws.Cells.Item(19,50).Value = 'Helena Thau'

# --- Row 20 (was row 19) ---
$ws.Cells.Item(20,1).Value = 112392981
$ws.Cells.Item(20,2).Value = 56430
$ws.Cells.Item(20,3).Value = 'Ovaliderad'
$ws.Cells.Item(20,4).Value = 'NT'
$ws.Cells.Item(20,5).Value = 100109
$ws.Cells.Item(20,6).Value = 'Tretåig hackspett'
$ws.Cells.Item(20,7).Value = 'Picoides tridactylus'
$ws.Cells.Item(20,8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(20,13).Value = 'färska spår'
$ws.Cells.Item(20,16).Value = 'Kärmsjöbäckens naturreservat (Kärmsjöbäckens naturreservat), Ång'
$ws.Cells.Item(20,17).Value = 583026
$ws.Cells.Item(20,18).Value = 7086695
$ws.Cells.Item(20,19).Value = 20
$ws.Cells.Item(20,20).Value = 'Västernorrland'
$ws.Cells.Item(20,21).Value = 'Sollefteå'
$ws.Cells.Item(20,22).Value = 'Ångermanland'
$ws.Cells.Item(20,23).Value = 'Junsele'
$ws.Cells.Item(20,25).NumberFormat = "@"
$ws.Cells.Item(20,25).Value = '2023-09-29'
$ws.Cells.Item(20,27).NumberFormat = "@"
$ws.Cells.Item(20,27).Value = '2023-09-29'
$ws.Cells.Item(20,30).Value = $false
$ws.Cells.Item(20,31).Value = $false
$ws.Cells.Item(20,33).Value = $false
$ws.Cells.Item(20,49).Value = 'Daniel Rutschman'
$ws.Cells.Item(20,50).Value = 'Daniel Rutschman'

# --- Row 21 (was row 20) ---
$ws.Cells.Item(21,1).Value = 112394976
$ws.Cells.Item(21,2).Value = 89553
$ws.Cells.Item(21,3).Value = 'Ovaliderad'
$ws.Cells.Item(21,4).Value = 'NT'
$ws.Cells.Item(21,5).Value = 1202
$ws.Cells.Item(21,6).Value = 'Ullticka'
$ws.Cells.Item(21,7).Value = 'Phellinidium ferrugineofuscum'
$ws.Cells.Item(21,8).Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Cells.Item(21,13).ClearContents()
$ws.Cells.Item(21,16).Value = 'Kärmsjöbäckens naturreservat (Kärmsjöbäckens naturreservat), Ång'
$ws.Cells.Item(21,17).Value = 583147
$ws.Cells.Item(21,18).Value = 7086540
$ws.Cells.Item(21,19).Value = 20
$ws.Cells.Item(21,20).Value = 'Västernorrland'
$ws.Cells.Item(21,21).Value = 'Sollefteå'
$ws.Cells.Item(21,22).Value = 'Ångermanland'
$ws.Cells.Item(21,23).Value = 'Junsele'
$ws.Cells.Item(21,25).NumberFormat = "@"
$ws.Cells.Item(21,25).Value = '2023-09-29'
$ws.Cells.Item(21,26).ClearContents()
$ws.Cells.Item(21,27).NumberFormat = "@"
$ws.Cells.Item(21,27).Value = '2023-09-29'
$ws.Cells.Item(21,28).ClearContents()
$ws.Cells.Item(21,30).Value = $false
$ws.Cells.Item(21,31).Value = $false
$ws.Cells.Item(21,33).Value = $false
$ws.Cells.Item(21,49).Value = 'Daniel Rutschman'
$ws.Cells.Item(21,50).Value = 'Daniel Rutschman'

# --- Row 23 (was row 28) ---
$ws.Cells.Item(23,1).Value = 112393743
$ws.Cells.Item(23,2).Value = 56430
$ws.Cells.Item(23,3).Value = 'Ovaliderad'
$ws.Cells.Item(23,4).Value = 'NT'
$ws.Cells.Item(23,5).Value = 100109
$ws.Cells.Item(23,6).Value = 'Tretåig hackspett'
$ws.Cells.Item(23,7).Value = 'Picoides tridactylus'
$ws.Cells.Item(23,8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(23,13).Value = 'äldre spår'
$ws.Cells.Item(23,16).Value = 'Stor-Kärmsjön, Ång'
$ws.Cells.Item(23,17).Value = 583023
$ws.Cells.Item(23,18).Value = 7086587
$ws.Cells.Item(23,19).Value = 20
$ws.Cells.Item(23,20).Value = 'Västernorrland'
$ws.Cells.Item(23,21).Value = 'Sollefteå'
$ws.Cells.Item(23,22).Value = 'Ångermanland'
$ws.Cells.Item(23,23).Value = 'Junsele'
$ws.Cells.Item(23,25).NumberFormat = "@"
$ws.Cells.Item(23,25).Value = '2023-09-29'
$ws.Cells.Item(23,26).ClearContents()
$ws.Cells.Item(23,27).NumberFormat = "@"
$ws.Cells.Item(23,27).Value = '2023-09-29'
$ws.Cells.Item(23,28).ClearContents()
$ws.Cells.Item(23,30).Value = $false
$ws.Cells.Item(23,31).Value = $false
$ws.Cells.Item(23,33).Value = $false
$ws.Cells.Item(23,49).Value = 'Michaela Ehmann'
$ws.Cells.Item(23,50).Value = 'Michaela Ehmann'

# --- Row 28 (was row 23) ---
$ws.Cells.Item(28,1).Value = 112392979
$ws.Cells.Item(28,2).Value = 77650
$ws.Cells.Item(28,3).Value = 'Ovaliderad'
$ws.Cells.Item(28,4).Value = 'NT'
$ws.Cells.Item(28,5).Value = 6425
$ws.Cells.Item(28,6).Value = 'Garnlav'
$ws.Cells.Item(28,7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(28,8).Value = '(Ach.) Ach.'
$ws.Cells.Item(28,13).ClearContents()
$ws.Cells.Item(28,16).Value = 'Kärmsjöbäckens naturreservat (Kärmsjöbäckens naturreservat), Ång'
$ws.Cells.Item(28,17).Value = 583017
$ws.Cells.Item(28,18).Value = 7086682
$ws.Cells.Item(28,19).Value = 10
$ws.Cells.Item(28,20).Value = 'Västernorrland'
$ws.Cells.Item(28,21).Value = 'Sollefteå'
$ws.Cells.Item(28,22).Value = 'Ångermanland'
$ws.Cells.Item(28,23).Value = 'Junsele'
$ws.Cells.Item(28,25).NumberFormat = "@"
$ws.Cells.Item(28,25).Value = '2023-09-29'
$ws.Cells.Item(28,26).NumberFormat = "@"
$ws.Cells.Item(28,26).Value = '14:10'
$ws.Cells.Item(28,27).NumberFormat = "@"
$ws.Cells.Item(28,27).Value = '2023-09-29'
$ws.Cells.Item(28,28).NumberFormat = "@"
$ws.Cells.Item(28,28).Value = '14:10'
$ws.Cells.Item(28,30).Value = $false
$ws.Cells.Item(28,31).Value = $false
$ws.Cells.Item(28,33).Value = $false
$ws.Cells.Item(28,49).Value = 'Helena Thau'
$ws.Cells.Item(28,50).Value = 'Helena Thau'

# --- Row 35 (was row 37) ---
$ws.Cells.Item(35,1).Value = 112394482
$ws.Cells.Item(35,2).Value = 90814
$ws.Cells.Item(35,3).Value = 'Ovaliderad'
$ws.Cells.Item(35,4).Value = 'LC'
$ws.Cells.Item(35,5).Value = 4364
$ws.Cells.Item(35,6).Value = 'Dropptaggsvamp'
$ws.Cells.Item(35,7).Value = 'Hydnellum ferrugineum'
$ws.Cells.Item(35,8).Value = '(Fr.:Fr.) P. Karst.'
$ws.Cells.Item(35,13).ClearContents()
$ws.Cells.Item(35,16).Value = 'Kärmsjöbäckens naturreservat (Kärmsjöbäckens naturreservat), Ång'
$ws.Cells.Item(35,17).Value = 583009
$ws.Cells.Item(35,18).Value = 7086522
$ws.Cells.Item(35,19).Value = 20
$ws.Cells.Item(35,20).Value = 'Västernorrland'
$ws.Cells.Item(35,21).Value = 'Sollefteå'
$ws.Cells.Item(35,22).Value = 'Ångermanland'
$ws.Cells.Item(35,23).Value = 'Junsele'
$ws.Cells.Item(35,25).NumberFormat = "@"
$ws.Cells.Item(35,25).Value = '2023-08-20'
$ws.Cells.Item(35,26).ClearContents()
$ws.Cells.Item(35,27).NumberFormat = "@"
$ws.Cells.Item(35,27).Value = '2023-08-20'
$ws.Cells.Item(35,28).ClearContents()
$ws.Cells.Item(35,30).Value = $false
$ws.Cells.Item(35,31).Value = $false
$ws.Cells.Item(35,33).Value = $false
$ws.Cells.Item(35,49).Value = 'Daniel Rutschman'
$ws.Cells.Item(35,50).Value = 'Daniel Rutschman'

# --- Row 36 (was row 38) ---
$ws.Cells.Item(36,1).Value = 112410306
$ws.Cells.Item(36,2).Value = 88637
$ws.Cells.Item(36,3).Value = 'Ovaliderad'
$ws.Cells.Item(36,4).Value = 'NT'
$ws.Cells.Item(36,5).Value = 1962
$ws.Cells.Item(36,6).Value = 'Vaddporing'
$ws.Cells.Item(36,7).Value = 'Anomoporia kamtschatica'
$ws.Cells.Item(36,8).Value = '(Parmasto) Bondartseva'
$ws.Cells.Item(36,16).Value = 'Kärmsjöbäckens naturreservat, Ång'
$ws.Cells.Item(36,17).Value = 583095
$ws.Cells.Item(36,18).Value = 7086753
$ws.Cells.Item(36,19).Value = 20
$ws.Cells.Item(36,20).Value = 'Västernorrland'
$ws.Cells.Item(36,21).Value = 'Sollefteå'
$ws.Cells.Item(36,22).Value = 'Ångermanland'
$ws.Cells.Item(36,23).Value = 'Junsele'
$ws.Cells.Item(36,25).NumberFormat = "@"
$ws.Cells.Item(36,25).Value = '2023-09-29'
$ws.Cells.Item(36,26).ClearContents()
$ws.Cells.Item(36,27).NumberFormat = "@"
$ws.Cells.Item(36,27).Value = '2023-09-29'
$ws.Cells.Item(36,28).ClearContents()
$ws.Cells.Item(36,30).Value = $false
$ws.Cells.Item(36,31).Value = $false
$ws.Cells.Item(36,33).Value = $false
$ws.Cells.Item(36,49).Value = 'Daniel Rutschman'
$ws.Cells.Item(36,50).Value = 'Daniel Rutschman'

# --- Row 37 (was row 35) ---
$ws.Cells.Item(37,1).Value = 112390945
$ws.Cells.Item(37,2).Value = 56430
$ws.Cells.Item(37,3).Value = 'Ovaliderad'
$ws.Cells.Item(37,4).Value = 'NT'
$ws.Cells.Item(37,5).Value = 100109
$ws.Cells.Item(37,6).Value = 'Tretåig hackspett'
$ws.Cells.Item(37,7).Value = 'Picoides tridactylus'
$ws.Cells.Item(37,8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(37,13).Value = 'färska spår'
$ws.Cells.Item(37,16).Value = 'Kärmsjöbäckens naturreservat (Kärmsjöbäckens naturreservat), Ång'
$ws.Cells.Item(37,17).Value = 583104
$ws.Cells.Item(37,18).Value = 7086760
$ws.Cells.Item(37,19).Value = 25
$ws.Cells.Item(37,20).Value = 'Västernorrland'
$ws.Cells.Item(37,21).Value = 'Sollefteå'
$ws.Cells.Item(37,22).Value = 'Ångermanland'
$ws.Cells.Item(37,23).Value = 'Junsele'
$ws.Cells.Item(37,25).NumberFormat = "@"
$ws.Cells.Item(37,25).Value = '2023-09-29'
$ws.Cells.Item(37,26).NumberFormat = "@"
$ws.Cells.Item(37,26).Value = '12:01'
$ws.Cells.Item(37,27).NumberFormat = "@"
$ws.Cells.Item(37,27).Value = '2023-09-29'
$ws.Cells.Item(37,28).NumberFormat = "@"
$ws.Cells.Item(37,28).Value = '12:01'
$ws.Cells.Item(37,30).Value = $false
$ws.Cells.Item(37,31).Value = $false
$ws.Cells.Item(37,33).Value = $false
$ws.Cells.Item(37,49).Value = 'Helena Thau'
$ws.Cells.Item(37,50).Value = 'Helena Thau'

# --- Row 38 (was row 36) ---
$ws.Cells.Item(38,1).Value = 112394782
$ws.Cells.Item(38,2).Value = 77402
$ws.Cells.Item(38,3).Value = 'Ovaliderad'
$ws.Cells.Item(38,4).Value = 'NT'
$ws.Cells.Item(38,5).Value = 6446
$ws.Cells.Item(38,6).Value = 'Kolflarnlav'
$ws.Cells.Item(38,7).Value = 'Carbonicola anthracophila'
$ws.Cells.Item(38,8).Value = '(Nyl.) Bendiksby & Timdal'
$ws.Cells.Item(38,16).Value = 'Kärmsjöbäckens naturreservat (Kärmsjöbäckens naturreservat), Ång'
$ws.Cells.Item(38,17).Value = 583127
$ws.Cells.Item(38,18).Value = 7086464
$ws.Cells.Item(38,19).Value = 10
$ws.Cells.Item(38,20).Value = 'Västernorrland'
$ws.Cells.Item(38,21).Value = 'Sollefteå'
$ws.Cells.Item(38,22).Value = 'Ångermanland'
$ws.Cells.Item(38,23).Value = 'Junsele'
$ws.Cells.Item(38,25).NumberFormat = "@"
$ws.Cells.Item(38,25).Value = '2023-09-29'
$ws.Cells.Item(38,26).NumberFormat = "@"
$ws.Cells.Item(38,26).Value = '14:10'
$ws.Cells.Item(38,27).NumberFormat = "@"
$ws.Cells.Item(38,27).Value = '2023-09-29'
$ws.Cells.Item(38,28).NumberFormat = "@"
$ws.Cells.Item(38,28).Value = '14:10'
$ws.Cells.Item(38,30).Value = $false
$ws.Cells.Item(38,31).Value = $false
$ws.Cells.Item(38,33).Value = $false
$ws.Cells.Item(38,49).Value = 'Helena Thau'
$ws.Cells.Item(38,50).Value = 'Helena Thau'
